$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K -> new F:M), shifting existing data right
$ws.Columns("D:E").Insert()

# Copy cell formatting (number format/style) from the (now-shifted) F:G columns into the
# new D:E columns, for exactly the rows that carry data (skip header-only rows with no D:K cells)
$fmtSrc = $excel.Union($ws.Range("F7:G35"), $ws.Range("F38:G77"), $ws.Range("F80:G102"))
$fmtDst = $excel.Union($ws.Range("D7:E35"), $ws.Range("D38:E77"), $ws.Range("D80:E102"))
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D (and E) column values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 5190000
$ws.Range("E8").Value = 5582000
$ws.Range("D9").Value = 313000
$ws.Range("E9").Value = 315000
$ws.Range("D10").Value = 4877000
$ws.Range("E10").Value = 5267000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 64000
$ws.Range("E14").Value = 5000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 4410000
$ws.Range("E17").Value = 4383000
$ws.Range("D18").Value = 780000
$ws.Range("E18").Value = 1199000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 780000
$ws.Range("E21").Value = 1199000
$ws.Range("D22").Value = 58000
$ws.Range("E22").Value = 53000
$ws.Range("D23").Value = 722000
$ws.Range("E23").Value = 1146000
$ws.Range("D24").Value = 190000
$ws.Range("E24").Value = 290000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 532000
$ws.Range("E26").Value = 856000
$ws.Range("D27").Value = 532000
$ws.Range("E27").Value = 856000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -7000
$ws.Range("E29").Value = -11000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 525000
$ws.Range("E33").Value = 845000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 525000
$ws.Range("E35").Value = 845000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 4337000
$ws.Range("E41").Value = 3429000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 851000
$ws.Range("E43").Value = 784000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 121906000
$ws.Range("E47").Value = 120785000
$ws.Range("D48").Value = 443000
$ws.Range("E48").Value = 435000
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 140406000
$ws.Range("E54").Value = 137941000
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 17308000
$ws.Range("E59").Value = 18055000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 5778000
$ws.Range("E61").Value = 5279000
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 116944000
$ws.Range("E66").Value = 114707000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 31788000
$ws.Range("E72").Value = 31461000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 23462000
$ws.Range("E76").Value = 23234000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 525000
$ws.Range("E81").Value = 845000
$ws.Range("D83").Value = 313000
$ws.Range("E83").Value = 315000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1355000
$ws.Range("E89").Value = 1852000
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -410000
$ws.Range("E94").Value = -1729000
$ws.Range("D96").Value = -198000
$ws.Range("E96").Value = -199000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -103000
$ws.Range("E100").Value = -525000
$ws.Range("D101").Value = 66000
$ws.Range("E101").Value = -16000
$ws.Range("D102").Value = 908000
$ws.Range("E102").Value = -418000
